$d = $word.ActiveDocument

# The edit reverts a previous change: remove the trailing empty paragraph,
# the "Test" paragraph, and the "345" paragraph, leaving only the first
# paragraph ("2nd update") before the section properties.
$firstPara = $d.Paragraphs.Item(1)
$startPos = $firstPara.Range.End
$endPos = $d.Content.End

$r = $d.Range($startPos, $endPos)
$r.Delete()
